# pTHg_OutflowR.xlsx update
# - Rename the existing sheet to "4_pTHg_OutflowR_31ct" (the full/31-count series)
# - Duplicate it to a new "4_pTHg_OutflowR_23ct" sheet (the trimmed/23-count series)
#   by removing the 8-row block (rows 21-28) that belongs only to the 31ct series
# - Highlight that 8-row block on the 31ct sheet with a light-cyan fill
# - Re-point the absPath / revision metadata, fix up selections, and leave the
#   new 23ct sheet as the active tab, matching the authored workbook.

$wb = $excel.ActiveWorkbook

# ---- Sheet 1: rename to the "31ct" (31-count) series -----------------------
$ws1 = $wb.Worksheets.Item(1)
$ws1.Name = "4_pTHg_OutflowR_31ct"

# ---- Create the "23ct" series by duplicating sheet 1, then trimming it -----
$ws1.Copy($null, $ws1) | Out-Null
$ws2 = $wb.Worksheets.Item(2)
$ws2.Name = "4_pTHg_OutflowR_23ct"

# Remove the 8-row block (rows 21:28) unique to the 31ct series so the
# remaining rows (previously 29:33) shift up to 21:25.
$ws2.Rows("21:28").Delete()

# Tab color for the new sheet (light cyan highlight)
$ws2.Tab.Color = 16644525

# New sheet's own view/selection
$ws2.Range("E25").Select()

# ---- Highlight rows 21:28 on the 31ct sheet with the same light cyan -------
$ws1.Range("A21:D28").Interior.Color = 16644525

# Restore sheet1's remembered selection before switching away from it
$ws1.Range("H35").Select()

# ---- Make the 23ct sheet the active/front tab ------------------------------
$ws2.Activate()

Write-Output "done"
